$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Ephb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.071327
$ws.Range("H2").Value = 0.213981
$ws.Range("I2").Value = 0.03356605248408491
$ws.Range("J2").Value = 0.03356605248408491
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.150777
$ws.Range("N2").Value = 0.452331
$ws.Range("O2").Value = 0.009673301965868179
$ws.Range("P2").Value = 0.009673301965868179
$ws.Range("Q2").Value = 0.010754471079
$ws.Range("R2").Value = 0.096790239711
$ws.Range("S2").Value = 0.000324694561480733
$ws.Range("T2").Value = 0.000324694561480733

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Ephb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.071327
$ws.Range("H3").Value = 0.213981
$ws.Range("I3").Value = 0.03356605248408491
$ws.Range("J3").Value = 0.03356605248408491
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 12.30706533333333
$ws.Range("N3").Value = 36.921196
$ws.Range("O3").Value = 0.7895763895222844
$ws.Range("P3").Value = 0.7895763895222843
$ws.Range("Q3").Value = 0.8778260490306667
$ws.Range("R3").Value = 7.900434441276
$ws.Range("S3").Value = 0.02650296253089927
$ws.Range("T3").Value = 0.02650296253089926

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb3"
$ws.Range("C4").Value = "Ephb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.071327
$ws.Range("H4").Value = 0.213981
$ws.Range("I4").Value = 0.03356605248408491
$ws.Range("J4").Value = 0.03356605248408491
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.129079333333333
$ws.Range("N4").Value = 9.387238
$ws.Range("O4").Value = 0.2007503085118475
$ws.Range("P4").Value = 0.2007503085118475
$ws.Range("Q4").Value = 0.2231878416086667
$ws.Range("R4").Value = 2.008690574478
$ws.Range("S4").Value = 0.00673839539170491
$ws.Range("T4").Value = 0.00673839539170491

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb3"
$ws.Range("C5").Value = "Ephb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.308223
$ws.Range("H5").Value = 0.9246690000000001
$ws.Range("I5").Value = 0.1450478695977975
$ws.Range("J5").Value = 0.1450478695977975
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.150777
$ws.Range("N5").Value = 0.452331
$ws.Range("O5").Value = 0.009673301965868179
$ws.Range("P5").Value = 0.009673301965868179
$ws.Range("Q5").Value = 0.046472939271
$ws.Range("R5").Value = 0.418256453439
$ws.Range("S5").Value = 0.001403091842125366
$ws.Range("T5").Value = 0.001403091842125366

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb3"
$ws.Range("C6").Value = "Ephb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.308223
$ws.Range("H6").Value = 0.9246690000000001
$ws.Range("I6").Value = 0.1450478695977975
$ws.Range("J6").Value = 0.1450478695977975
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.30706533333333
$ws.Range("N6").Value = 36.921196
$ws.Range("O6").Value = 0.7895763895222844
$ws.Range("P6").Value = 0.7895763895222843
$ws.Range("Q6").Value = 3.793320598236
$ws.Range("R6").Value = 34.13988538412401
$ws.Range("S6").Value = 0.1145263731849281
$ws.Range("T6").Value = 0.1145263731849281

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb3"
$ws.Range("C7").Value = "Ephb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.308223
$ws.Range("H7").Value = 0.9246690000000001
$ws.Range("I7").Value = 0.1450478695977975
$ws.Range("J7").Value = 0.1450478695977975
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.129079333333333
$ws.Range("N7").Value = 9.387238
$ws.Range("O7").Value = 0.2007503085118475
$ws.Range("P7").Value = 0.2007503085118475
$ws.Range("Q7").Value = 0.964454219358
$ws.Range("R7").Value = 8.680087974222001
$ws.Range("S7").Value = 0.02911840457074407
$ws.Range("T7").Value = 0.02911840457074407

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efnb3"
$ws.Range("C8").Value = "Ephb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.745424333333333
$ws.Range("H8").Value = 5.236273
$ws.Range("I8").Value = 0.8213860779181176
$ws.Range("J8").Value = 0.8213860779181176
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.150777
$ws.Range("N8").Value = 0.452331
$ws.Range("O8").Value = 0.009673301965868179
$ws.Range("P8").Value = 0.009673301965868179
$ws.Range("Q8").Value = 0.263169844707
$ws.Range("R8").Value = 2.368528602363
$ws.Range("S8").Value = 0.007945515562262081
$ws.Range("T8").Value = 0.007945515562262081

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efnb3"
$ws.Range("C9").Value = "Ephb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.745424333333333
$ws.Range("H9").Value = 5.236273
$ws.Range("I9").Value = 0.8213860779181176
$ws.Range("J9").Value = 0.8213860779181176
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 12.30706533333333
$ws.Range("N9").Value = 36.921196
$ws.Range("O9").Value = 0.7895763895222844
$ws.Range("P9").Value = 0.7895763895222843
$ws.Range("Q9").Value = 21.48105130472311
$ws.Range("R9").Value = 193.329461742508
$ws.Range("S9").Value = 0.6485470538064571
$ws.Range("T9").Value = 0.6485470538064569

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efnb3"
$ws.Range("C10").Value = "Ephb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.745424333333333
$ws.Range("H10").Value = 5.236273
$ws.Range("I10").Value = 0.8213860779181176
$ws.Range("J10").Value = 0.8213860779181176
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.129079333333333
$ws.Range("N10").Value = 9.387238
$ws.Range("O10").Value = 0.2007503085118475
$ws.Range("P10").Value = 0.2007503085118475
$ws.Range("Q10").Value = 5.461571209330444
$ws.Range("R10").Value = 49.15414088397399
$ws.Range("S10").Value = 0.1648935085493985
$ws.Range("T10").Value = 0.1648935085493985
